# Append a new data row (row 3) to the portfolio worksheet, mirroring the
# structure/style of the existing row 2 (plain, unstyled data cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 3

$ws.Cells.Item($row, 1).Value  = 2                 # A3  Id
$ws.Cells.Item($row, 2).Value  = "ETF"              # B3  category
$ws.Cells.Item($row, 3).Value  = "MSC World"         # C3  assetName
$ws.Cells.Item($row, 4).Value  = "Fineco"            # D3  position
$ws.Cells.Item($row, 5).Value  = 3                 # E3  riskLevel
$ws.Cells.Item($row, 6).Value  = "sdr4558"           # F3  ticker
$ws.Cells.Item($row, 7).Value  = "ed38383838"        # G3  isin
$ws.Cells.Item($row, 8).Value  = "22/05/2000"        # H3  createdAt
$ws.Cells.Item($row, 9).Value  = 234                # I3  createdAmount
$ws.Cells.Item($row, 10).Value = 23.45              # J3  createdUnitPrice
$ws.Cells.Item($row, 11).Value = 5487.3             # K3  createdTotalValue
$ws.Cells.Item($row, 12).Value = "23/07/2024"        # L3  updatedAt
$ws.Cells.Item($row, 13).Value = 278                # M3  updatedAmount
$ws.Cells.Item($row, 14).Value = 25.89              # N3  updatedUnitPrice
$ws.Cells.Item($row, 15).Value = 7197.42            # O3  updatedTotalValue
$ws.Cells.Item($row, 16).Value = ""                 # P3  accumulationPlan (blank)
$ws.Cells.Item($row, 17).Value = 0                  # Q3  accumulationAmount
$ws.Cells.Item($row, 18).Value = 0                  # R3  incomePerYear
$ws.Cells.Item($row, 19).Value = 0                  # S3  rentalIncome
$ws.Cells.Item($row, 20).Value = "ETF di prova"      # T3  note
